# Apply the balance-ytd demo report edits:
#  - bump the report period end date in the summary sheet
#  - bump a set of document dates in the documents sheet (off-by-one day)
#  - change the currency number format (remove the fixed ₪ prefix, add red-negative variant)

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("סיכום")
$wsDocs = $wb.Worksheets.Item("מסמכים")

# --- Summary sheet: period label ---
$wsSummary.Range("A2").Value = "תקופה: 2026-01-01 עד 2026-02-13"

# --- Documents sheet: bump individual document dates by one day ---
$dateUpdates = @{
    "B10" = "2026-01-10"
    "B11" = "2026-01-11"
    "B18" = "2026-01-19"
    "B19" = "2026-01-20"
    "B20" = "2026-01-21"
    "B25" = "2026-01-27"
    "B26" = "2026-01-28"
    "B27" = "2026-01-29"
    "B28" = "2026-01-30"
    "B29" = "2026-01-31"
    "B30" = "2026-02-01"
    "B33" = "2026-02-05"
    "B34" = "2026-02-06"
    "B35" = "2026-02-07"
    "B36" = "2026-02-08"
    "B37" = "2026-02-09"
    "B38" = "2026-02-10"
    "B39" = "2026-02-11"
}

foreach ($addr in $dateUpdates.Keys) {
    $wsDocs.Range($addr).Value = $dateUpdates[$addr]
}

# --- Styles: update the ILS currency number format code ---
# numFmtId 164 goes from "₪#,##0.00" to "#,##0.00;[Red]-#,##0.00"
$wb.Styles.Item("ILS").NumberFormat = "#,##0.00;[Red]-#,##0.00"
